$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.927.25'
$ws.Range("E2").Value = '  -3.86%  '
$ws.Range("D3").Value = '3.021.44'
$ws.Range("E3").Value = '  -3.92%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '527.59'
$ws.Range("E5").Value = '  -6.08%  '
$ws.Range("D6").Value = '128.90'
$ws.Range("E6").Value = '  -9.27%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").Value = '3.013.51'
$ws.Range("E8").Value = '  -4.02%  '
$ws.Range("E9").Value = '  -1.73%  '
$ws.Range("E10").Value = '  -3.62%  '
$ws.Range("D11").Value = '5.99'
$ws.Range("E11").Value = '  -11.40%  '
$ws.Range("D12").Value = '0.443'
$ws.Range("E12").Value = '  -4.57%  '
$ws.Range("E13").Value = '  -1.18%  '
$ws.Range("D14").Value = '33.19'
$ws.Range("E14").Value = '  -9.26%  '
$ws.Range("D15").Value = '3.479.34'
$ws.Range("E15").Value = '  -4.56%  '
$ws.Range("D16").Value = '61.992.68'
$ws.Range("E16").Value = '  -3.87%  '
$ws.Range("D18").Value = '3.023.31'
$ws.Range("E18").Value = '  -3.98%  '
$ws.Range("D19").Value = '6.40'
$ws.Range("E19").Value = '  -5.97%  '
$ws.Range("D20").Value = '469.23'
$ws.Range("E20").Value = '  -9.02%  '
$ws.Range("D21").Value = '12.96'
$ws.Range("E21").Value = '  -7.23%  '
$ws.Range("D22").Value = '0.682'
$ws.Range("E22").Value = '  -4.47%  '
$ws.Range("D23").Value = '6.92'
$ws.Range("E23").Value = '  -6.76%  '
$ws.Range("D24").Value = '77.45'
$ws.Range("E24").Value = '  -1.77%  '
$ws.Range("D25").Value = '11.66'
$ws.Range("E25").Value = '  -8.28%  '
$ws.Range("B27").Value = 'RenderToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D27").Value = '8.00'
$ws.Range("E27").Value = '  -9.75%  '
$ws.Range("B28").Value = 'PancakeSwap'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D28").Value = '2.61'
$ws.Range("E28").Value = '  -7.49%  '
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("D30").Value = '25.21'
$ws.Range("E30").Value = '  -4.93%  '
$ws.Range("E31").Value = '  -15.60%  '
$ws.Range("D32").Value = '1.08'
$ws.Range("E32").Value = '  -4.48%  '
$ws.Range("E33").Value = '  -10.02%  '
$ws.Range("D34").Value = '55.88'
$ws.Range("E34").Value = '  +3.74%  '
$ws.Range("D35").Value = '5.09'
$ws.Range("E35").Value = '  -5.16%  '
$ws.Range("D36").Value = '5.74'
$ws.Range("E36").Value = '  -5.29%  '
$ws.Range("D37").Value = '457.82'
$ws.Range("E37").Value = '  -16.61%  '
$ws.Range("D38").Value = '3.037.86'
$ws.Range("E38").Value = '  -3.77%  '
$ws.Range("D39").Value = '0.0382'
$ws.Range("E39").Value = '  -11.73%  '
$ws.Range("D40").Value = '0.0766'
$ws.Range("E40").Value = '  -6.70%  '
$ws.Range("E41").Value = '  -8.94%  '
$ws.Range("D42").Value = '7.85'
$ws.Range("E42").Value = '  -4.88%  '
$ws.Range("D43").Value = '2.48'
$ws.Range("E43").Value = '  -9.27%  '
$ws.Range("E44").Value = '  +0.00%  '
$ws.Range("E45").Value = '  -8.36%  '
$ws.Range("D46").Value = '1.96'
$ws.Range("E46").Value = '  -11.49%  '
$ws.Range("E47").Value = '  -1.26%  '
$ws.Range("D48").Value = '23.65'
$ws.Range("E48").Value = '  -6.08%  '
$ws.Range("E49").Value = '  -2.38%  '
$ws.Range("D50").Value = '114.65'
$ws.Range("E50").Value = '  -4.98%  '
$ws.Range("B51").Value = 'ThetaToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D51").Value = '1.91'
$ws.Range("E51").Value = '  -8.83%  '
